# Sync file from Google Drive
# Refreshes the NextBus arrival snapshot (EstimatedTimeOfArrival, Load,
# Monitored, TypeOfBus, MinutesToArrival) on each of the three NextBus
# sheets to the newly-pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NextBus1")
$ws.Range("F2").Value = 45690.5340162037
$ws.Range("O2").Value = 2
$ws.Range("F3").Value = 45690.53364583333
$ws.Range("O3").Value = 2
$ws.Range("F4").Value = 45690.53380787037
$ws.Range("O4").Value = 2
$ws.Range("F5").Value = 45690.53730324074
$ws.Range("F6").Value = 45690.53711805555
$ws.Range("O6").Value = 7
$ws.Range("F7").Value = 45690.5315625
$ws.Range("F8").Value = 45690.53837962963
$ws.Range("L8").Value = "DD"
$ws.Range("O8").Value = 9
$ws.Range("F9").Value = 45690.5337037037
$ws.Range("O9").Value = 2
$ws.Range("F10").Value = 45690.53709490741
$ws.Range("O10").Value = 7
$ws.Range("F11").Value = 45690.53917824074
$ws.Range("O11").Value = 10
$ws.Range("F12").Value = 45690.53928240741
$ws.Range("L12").Value = "SD"
$ws.Range("O12").Value = 10
$ws.Range("F13").Value = 45690.53582175926
$ws.Range("L13").Value = "SD"
$ws.Range("O13").Value = 5
$ws.Range("F14").Value = 45690.53400462963
$ws.Range("O14").Value = 2
$ws.Range("F15").Value = 45690.53524305556
$ws.Range("O15").Value = 4
$ws = $wb.Worksheets.Item("NextBus2")
$ws.Range("F2").Value = 45690.54305555556
$ws.Range("L2").Value = "SD"
$ws.Range("O2").Value = 15
$ws.Range("F3").Value = 45690.54038194445
$ws.Range("O3").Value = 11
$ws.Range("F4").Value = 45690.54719907408
$ws.Range("O4").Value = 21
$ws.Range("F5").Value = 45690.54694444445
$ws.Range("O5").Value = 21
$ws.Range("F6").Value = 45690.54377314815
$ws.Range("L6").Value = "SD"
$ws.Range("O6").Value = 16
$ws.Range("F7").Value = 45690.54136574074
$ws.Range("I7").Value = "SDA"
$ws.Range("J7").Value = 1
$ws.Range("O7").Value = 13
$ws.Range("F8").Value = 45690.54532407408
$ws.Range("O8").Value = 19
$ws.Range("F9").Value = 45690.54324074074
$ws.Range("O9").Value = 16
$ws.Range("F10").Value = 45690.54424768518
$ws.Range("O10").Value = 17
$ws.Range("F11").Value = 45690.54575231481
$ws.Range("J11").Value = 0
$ws.Range("O11").Value = 19
$ws.Range("F12").Value = 45690.54619212963
$ws.Range("O12").Value = 20
$ws.Range("F13").Value = 45690.54462962963
$ws.Range("L13").Value = "DD"
$ws.Range("O13").Value = 18
$ws.Range("F14").Value = 45690.5415162037
$ws.Range("O14").Value = 13
$ws.Range("F15").Value = 45690.5421412037
$ws.Range("O15").Value = 14
$ws = $wb.Worksheets.Item("NextBus3")
$ws.Range("F2").Value = 45690.55394675926
$ws.Range("O2").Value = 31
$ws.Range("F3").Value = 45690.54622685185
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 20
$ws.Range("F4").Value = 45690.55393518518
$ws.Range("O4").Value = 31
$ws.Range("F5").Value = 45690.55587962963
$ws.Range("J5").Value = 1
$ws.Range("O5").Value = 34
$ws.Range("F6").Value = 45690.55451388889
$ws.Range("L6").Value = "SD"
$ws.Range("O6").Value = 32
$ws.Range("F7").Value = 45690.54631944445
$ws.Range("O7").Value = 20
$ws.Range("F8").Value = 45690.5522337963
$ws.Range("L8").Value = "SD"
$ws.Range("O8").Value = 29
$ws.Range("F9").Value = 45690.55381944445
$ws.Range("L9").Value = "DD"
$ws.Range("O9").Value = 31
$ws.Range("F10").Value = 45690.553125
$ws.Range("O10").Value = 30
$ws.Range("F11").Value = 45690.55371527778
$ws.Range("O11").Value = 31
$ws.Range("F12").Value = 45690.55106481481
$ws.Range("O12").Value = 27
$ws.Range("F13").Value = 45690.55675925926
$ws.Range("L13").Value = "DD"
$ws.Range("O13").Value = 35
$ws.Range("F14").Value = 45690.54907407407
$ws.Range("O14").Value = 24
$ws.Range("F15").Value = 45690.55275462963
$ws.Range("J15").Value = 1
$ws.Range("O15").Value = 29
